# Applies the "cryptos list" price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $range = $ws.Range($cellRef)
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        # Looks like a plain number (e.g. "0.999") - Excel would silently
        # coerce Value to a Double and drop the exact text (trailing zeros,
        # etc). Use the text-prefix trick to force String storage, then
        # restore the default "Normal" style so no quote-prefix formatting
        # lingers on the cell.
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-TextValue "D2" "66.238.41"
Set-TextValue "E2" "  +7.56%  "
Set-TextValue "D3" "3.018.07"
Set-TextValue "E3" "  +4.68%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "583.49"
Set-TextValue "E5" "  +3.11%  "
Set-TextValue "D6" "156.70"
Set-TextValue "E6" "  +9.98%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.11%  "
Set-TextValue "D8" "3.014.62"
Set-TextValue "E8" "  +4.61%  "
Set-TextValue "D10" "6.98"
Set-TextValue "E10" "  +1.80%  "
Set-TextValue "E11" "  +7.11%  "
Set-TextValue "E12" "  +5.69%  "
Set-TextValue "D13" "0.0000252"
Set-TextValue "E13" "  +10.00%  "
Set-TextValue "D14" "34.53"
Set-TextValue "E14" "  +9.48%  "
Set-TextValue "E15" "  +0.70%  "
Set-TextValue "D16" "66.163.95"
Set-TextValue "E16" "  +7.46%  "
Set-TextValue "D17" "3.515.46"
Set-TextValue "E17" "  +4.53%  "
Set-TextValue "E18" "  +7.10%  "
Set-TextValue "D19" "3.020.39"
Set-TextValue "E19" "  +4.62%  "
Set-TextValue "D20" "463.39"
Set-TextValue "E20" "  +8.15%  "
Set-TextValue "D21" "13.91"
Set-TextValue "E21" "  +7.17%  "
Set-TextValue "D22" "0.684"
Set-TextValue "E22" "  +5.32%  "
Set-TextValue "E23" "  +8.48%  "
Set-TextValue "D24" "82.24"
Set-TextValue "E24" "  +4.33%  "
Set-TextValue "D25" "2.26"
Set-TextValue "E25" "  +13.14%  "
Set-TextValue "D26" "12.47"
Set-TextValue "E26" "  +5.43%  "
Set-TextValue "D27" "10.65"
Set-TextValue "E27" "  +7.64%  "
Set-TextValue "E28" "  -0.01%  "
Set-TextValue "D29" "8.03"
Set-TextValue "E29" "  +14.33%  "
Set-TextValue "D30" "2.38"
Set-TextValue "E30" "  +17.86%  "
Set-TextValue "E31" "  +0.60%  "
Set-TextValue "E32" "  +5.06%  "
Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  -0.08%  "
Set-TextValue "D36" "0.995"
Set-TextValue "E36" "  +4.28%  "
Set-TextValue "D37" "5.78"
Set-TextValue "E37" "  +8.23%  "
Set-TextValue "E38" "  +14.39%  "
Set-TextValue "E39" "  +9.49%  "
Set-TextValue "D40" "49.52"
Set-TextValue "E40" "  +1.53%  "
Set-TextValue "D41" "0.122"
Set-TextValue "E41" "  +7.72%  "
Set-TextValue "D42" "0.302"
Set-TextValue "E42" "  +13.88%  "
Set-TextValue "D43" "43.73"
Set-TextValue "E43" "  +11.17%  "
Set-TextValue "E44" "  +3.52%  "
Set-TextValue "D45" "391.13"
Set-TextValue "E45" "  +14.35%  "
Set-TextValue "D46" "2.805.21"
Set-TextValue "E46" "  +4.83%  "
Set-TextValue "E47" "  +6.24%  "
Set-TextValue "D48" "133.94"
Set-TextValue "E48" "  +1.10%  "
Set-TextValue "E49" "  -0.04%  "
Set-TextValue "D50" "23.57"
Set-TextValue "E50" "  +10.24%  "
Set-TextValue "D51" "0.107"
Set-TextValue "E51" "  +4.62%  "

# Row 33 / Row 34: re-sorted coins (EthereumClassic now ranks above Hedera),
# each with refreshed price and volume figures.
Set-TextValue "B33" "EthereumClassic"
Set-TextValue "C33" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "27.07"
Set-TextValue "E33" "  +6.82%  "

Set-TextValue "B34" "Hedera"
Set-TextValue "C34" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D34" "0.111"
Set-TextValue "E34" "  +5.13%  "
